# Apply the commit's numeric answer updates to the multiplication table.
$d = $word.ActiveDocument

$replacements = @(
    @("22×66=1452", "82×74=6068"),
    @("82×98=8036", "43×38=1634"),
    @("38×88=3344", "54×34=1836"),
    @("30×41=1230", "86×93=7998"),
    @("73×40=2920", "85×50=4250"),
    @("16×87=1392", "35×50=1750"),
    @("98×37=3626", "86×14=1204"),
    @("57×98=5586", "95×55=5225"),
    @("59×25=1475", "25×67=1675"),
    @("86×82=7052", "69×35=2415"),
    @("94×94=8836", "12×39=468"),
    @("83×74=6142", "28×90=2520"),
    @("11×31=341",  "70×49=3430"),
    @("12×57=684",  "15×54=810"),
    @("71×66=4686", "21×56=1176"),
    @("11×46=506",  "18×92=1656"),
    @("68×19=1292", "25×96=2400"),
    @("97×39=3783", "31×30=930"),
    @("39×56=2184", "47×62=2914"),
    @("13×96=1248", "84×65=5460"),
    @("42×80=3360", "76×82=6232"),
    @("87×66=5742", "33×55=1815"),
    @("67×43=2881", "26×19=494"),
    @("92×40=3680", "41×38=1558"),
    @("74×40=2960", "48×20=960")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
